# Update odds values in Sheet1 as per the FlashScore 2024-10-16 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 5 ---
$ws.Range("N5").Value = 8

# --- Row 6 ---
$ws.Range("I6").Value = 1.57

# --- Row 7 ---
$ws.Range("G7").Value = 1.67

# --- Row 8 ---
$ws.Range("M8").Value = 8.300000000000001

# --- Row 14 ---
$ws.Range("G14").Value = 2.77
$ws.Range("I14").Value = 2.32
$ws.Range("J14").Value = 3.35
$ws.Range("K14").Value = 2.1
$ws.Range("L14").Value = 2.95
$ws.Range("N14").Value = 6.9
$ws.Range("O14").Value = 1.34
$ws.Range("P14").Value = 3
$ws.Range("Q14").Value = 2.02
$ws.Range("W14").Value = 8.5
$ws.Range("X14").Value = 14
$ws.Range("Y14").Value = 10.25
$ws.Range("Z14").Value = 32
$ws.Range("AA14").Value = 24
$ws.Range("AB14").Value = 35
$ws.Range("AC14").Value = 6.9
$ws.Range("AE14").Value = 14.5
$ws.Range("AH14").Value = 7.6
$ws.Range("AI14").Value = 11.25
$ws.Range("AJ14").Value = 9.25
$ws.Range("AK14").Value = 24
$ws.Range("AL14").Value = 20
$ws.Range("AN14").Value = 4.75
$ws.Range("AO14").Value = 15
$ws.Range("AP14").Value = 23
$ws.Range("AQ14").Value = 70
$ws.Range("AU14").Value = 7.1
$ws.Range("AV14").Value = 65
$ws.Range("AW14").Value = 4.3
$ws.Range("AX14").Value = 12.5
$ws.Range("AZ14").Value = 50
